$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off "Share" (column K) for the clone scenes so that a payer logging in
# with a single clone scene creates a new group instead of joining a shared one.
$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("K15").Value = 0

# Update the active selection to reflect where the edit was made.
$ws.Range("K10").Select()
